$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = "한국 언론의 ‘빅데이터’ 기사 수준"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/kr-press-big-data/#utm_source=rss&utm_medium=rss&utm_campaign=kr-press-big-data"

$ws.Range("D28").Value = "Mobile manipulator 101 :: Combined Jacobian & Dual trajectory"
$ws.Range("E28").Value = "https://ropiens.tistory.com/197"
